$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2339561.2
$ws.Range("I33").Value = 2620226
$ws.Range("K33").Value = 2620226
$ws.Range("M33").Value = -2619997

$ws.Range("H58").Value = 1750
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = $null

$ws.Range("H62").Value = 4033.25
$ws.Range("I62").Value = 2961.8572
$ws.Range("J62").Value = 4866.5557
$ws.Range("K62").Value = 2961.8572
$ws.Range("L62").Value = 4866.5557
$ws.Range("M62").Value = -2337.8572
$ws.Range("N62").Value = -6114.5557

$ws.Range("H65").Value = 4033.25
$ws.Range("I65").Value = 2961.8572
$ws.Range("J65").Value = 4866.5557
$ws.Range("K65").Value = 14809.286
$ws.Range("L65").Value = 24332.7785
$ws.Range("M65").Value = -11689.286
$ws.Range("N65").Value = -30572.7785

$ws.Range("H76").Value = 4999.9165
$ws.Range("I76").Value = 4999.857
$ws.Range("K76").Value = 4999.857
$ws.Range("M76").Value = -4684.857

$ws.Range("H79").Value = 4999.9165
$ws.Range("I79").Value = 4999.857
$ws.Range("K79").Value = 4999.857
$ws.Range("M79").Value = -3907.857

$ws.Range("H94").Value = 864.6667
$ws.Range("I94").Value = 864.6667
$ws.Range("K94").Value = 864.6667
$ws.Range("M94").Value = -413.6667

$ws.Range("H98").Value = 1756
$ws.Range("I98").Value = 786.1111
$ws.Range("K98").Value = 786.1111
$ws.Range("M98").Value = 711.8889

$ws.Range("H100").Value = 2703
$ws.Range("I100").Value = 1804.75
$ws.Range("J100").Value = 5098.3335
$ws.Range("K100").Value = 1804.75
$ws.Range("L100").Value = 5098.3335
$ws.Range("M100").Value = -1263.75
$ws.Range("N100").Value = -6180.3335

$ws.Range("H112").Value = 1219.862
$ws.Range("J112").Value = 1245.2693
$ws.Range("L112").Value = 3735.8079
$ws.Range("N112").Value = -5951.8079

$ws.Range("H116").Value = 102370.91
$ws.Range("I116").Value = 116754.89
$ws.Range("J116").Value = 91582.914
$ws.Range("K116").Value = 116754.89
$ws.Range("L116").Value = 91582.914
$ws.Range("M116").Value = -113312.89
$ws.Range("N116").Value = -98466.914

$ws.Range("H122").Value = 1756
$ws.Range("I122").Value = 786.1111
$ws.Range("K122").Value = 2358.3333
$ws.Range("M122").Value = 91.66670000000022

$ws.Range("H132").Value = 50784.81
$ws.Range("I132").Value = 79183.62
$ws.Range("K132").Value = 237550.86
$ws.Range("M132").Value = -235020.86

$ws.Range("H137").Value = 4006.6924
$ws.Range("I137").Value = 1856.7142
$ws.Range("J137").Value = 6515
$ws.Range("K137").Value = 5570.142599999999
$ws.Range("L137").Value = 19545
$ws.Range("M137").Value = -3020.142599999999
$ws.Range("N137").Value = -24645

$ws.Range("H138").Value = 3771.6072
$ws.Range("J138").Value = 5051
$ws.Range("L138").Value = 15153
$ws.Range("N138").Value = -25433

$ws.Range("H141").Value = 6934.778
$ws.Range("I141").Value = 6060.5713
$ws.Range("K141").Value = 18181.7139
$ws.Range("M141").Value = -13001.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5854.75
$ws.Range("I2").Value = 6377.1924
$ws.Range("K2").Value = 6377.1924
$ws.Range("M2").Value = -6264.1924

$ws.Range("H116").Value = 5854.75
$ws.Range("I116").Value = 6377.1924
$ws.Range("K116").Value = 6377.1924
$ws.Range("M116").Value = -4083.1924

$ws.Range("H122").Value = 5832.3335
$ws.Range("I122").Value = 4997.5
$ws.Range("J122").Value = 6249.75
$ws.Range("K122").Value = 14992.5
$ws.Range("L122").Value = 18749.25
$ws.Range("M122").Value = -12542.5
$ws.Range("N122").Value = -23649.25

$ws.Range("H132").Value = 14290534
$ws.Range("I132").Value = 2898.2
$ws.Range("J132").Value = 50009624
$ws.Range("K132").Value = 8694.599999999999
$ws.Range("L132").Value = 150028872
$ws.Range("M132").Value = -6164.599999999999
$ws.Range("N132").Value = -150033932

$ws.Range("H141").Value = 28000
$ws.Range("I141").Value = 28000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 28000
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -22820
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5854.75
$ws.Range("I3").Value = 6377.1924
$ws.Range("K3").Value = 6377.1924
$ws.Range("M3").Value = -6263.1924

$ws.Range("H20").Value = 1708.2778
$ws.Range("I20").Value = 1390.3846
$ws.Range("J20").Value = 2534.8
$ws.Range("K20").Value = 1390.3846
$ws.Range("L20").Value = 2534.8
$ws.Range("M20").Value = -1143.3846
$ws.Range("N20").Value = -3028.8

$ws.Range("H50").Value = 63333.332
$ws.Range("J50").Value = 75000
$ws.Range("L50").Value = 75000
$ws.Range("N50").Value = -76148

$ws.Range("H75").Value = 21400
$ws.Range("I75").Value = 21400
$ws.Range("K75").Value = 21400
$ws.Range("M75").Value = -20464

$ws.Range("H78").Value = 21400
$ws.Range("I78").Value = 21400
$ws.Range("K78").Value = 64200
$ws.Range("M78").Value = -59520

$ws.Range("H105").Value = 3175
$ws.Range("I105").Value = 3465.5715
$ws.Range("J105").Value = 2884.4285
$ws.Range("K105").Value = 3465.5715
$ws.Range("L105").Value = 2884.4285
$ws.Range("M105").Value = -1718.5715
$ws.Range("N105").Value = -6378.4285

$ws.Range("H107").Value = 26614.334
$ws.Range("I107").Value = 7154.316
$ws.Range("J107").Value = 100562.4
$ws.Range("K107").Value = 7154.316
$ws.Range("L107").Value = 100562.4
$ws.Range("M107").Value = -5234.316
$ws.Range("N107").Value = -104402.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4442.6665
$ws.Range("I86").Value = 3328
$ws.Range("K86").Value = 3328
$ws.Range("M86").Value = -2205

$ws.Range("H89").Value = 4442.6665
$ws.Range("I89").Value = 3328
$ws.Range("K89").Value = 16640
$ws.Range("M89").Value = -11024

$ws.Range("H94").Value = 622.82355
$ws.Range("I94").Value = 515.4
$ws.Range("K94").Value = 515.4
$ws.Range("M94").Value = -64.39999999999998

$ws.Range("H122").Value = 26630.1
$ws.Range("I122").Value = 1666.8334
$ws.Range("K122").Value = 5000.5002
$ws.Range("M122").Value = -2550.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 614.86957
$ws.Range("J2").Value = 202.13333
$ws.Range("L2").Value = 1212.79998
$ws.Range("N2").Value = -1438.79998

$ws.Range("H26").Value = 680.36365
$ws.Range("J26").Value = 698.2222
$ws.Range("L26").Value = 2094.6666
$ws.Range("N26").Value = -2670.6666

$ws.Range("H130").Value = 3999.5
$ws.Range("I130").Value = 3999.5
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 11998.5
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -6978.5
$ws.Range("N130").Value = $null

$ws.Range("H131").Value = 569416.6
$ws.Range("I131").Value = 864.9091
$ws.Range("K131").Value = 2594.7273
$ws.Range("M131").Value = 2445.2727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = $null

$ws.Range("H70").Value = 6351.0415
$ws.Range("I70").Value = 6344
$ws.Range("K70").Value = 6344
$ws.Range("M70").Value = -6074

$ws.Range("H73").Value = 6351.0415
$ws.Range("I73").Value = 6344
$ws.Range("K73").Value = 6344
$ws.Range("M73").Value = -5408

$ws.Range("H97").Value = 1660.6666
$ws.Range("I97").Value = 991
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 991
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -495
$ws.Range("N97").Value = -3992

$ws.Range("H102").Value = 3562.1191
$ws.Range("I102").Value = 3416.1333
$ws.Range("K102").Value = 3416.1333
$ws.Range("M102").Value = -1794.1333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1250.3334
$ws.Range("I93").Value = 698
$ws.Range("J93").Value = 1526.5
$ws.Range("K93").Value = 698
$ws.Range("L93").Value = 1526.5
$ws.Range("M93").Value = 550
$ws.Range("N93").Value = -4022.5

$ws.Range("H132").Value = 3874.25
$ws.Range("I132").Value = 3999
$ws.Range("J132").Value = 3749.5
$ws.Range("K132").Value = 11997
$ws.Range("L132").Value = 11248.5
$ws.Range("M132").Value = -9467
$ws.Range("N132").Value = -16308.5

$ws.Range("H141").Value = 126749
$ws.Range("J141").Value = 126749
$ws.Range("L141").Value = 126749
$ws.Range("N141").Value = -137109

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 23000
$ws.Range("I18").Value = 23000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 23000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -22827
$ws.Range("N18").Value = $null

$ws.Range("H74").Value = 163871.8
$ws.Range("J74").Value = 163871.8
$ws.Range("L74").Value = 163871.8
$ws.Range("N74").Value = -165743.8

$ws.Range("H77").Value = 163871.8
$ws.Range("J77").Value = 163871.8
$ws.Range("L77").Value = 491615.4
$ws.Range("N77").Value = -500975.4

$ws.Range("H81").Value = 7629.75
$ws.Range("I81").Value = 9673.333000000001
$ws.Range("K81").Value = 19346.666
$ws.Range("M81").Value = -18285.666

$ws.Range("H84").Value = 7629.75
$ws.Range("I84").Value = 9673.333000000001
$ws.Range("K84").Value = 96733.33
$ws.Range("M84").Value = -91429.33

$ws.Range("H136").Value = 1560.0714
$ws.Range("I136").Value = 1394.6364
$ws.Range("J136").Value = 2166.6667
$ws.Range("K136").Value = 4183.9092
$ws.Range("L136").Value = 6500.000100000001
$ws.Range("M136").Value = -1633.9092
$ws.Range("N136").Value = -11600.0001
